$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "4x4 Squat Racks": two price updates
# ---------------------------------------------------------------------------
$wsRacks = $wb.Worksheets.Item("4x4 Squat Racks")

# C2: $2,149.00 -> $2,155.00 (keep as plain text, not an auto-converted number)
$wsRacks.Range("C2").Value = "'$2,155.00"
$wsRacks.Range("C2").Style = "Normal"

# C4: $1,520.00 -> Price not available
$wsRacks.Range("C4").Value = "Price not available"
$wsRacks.Range("C4").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "Squat Stands": one price update + a brand new product row
# ---------------------------------------------------------------------------
$wsStands = $wb.Worksheets.Item("Squat Stands")

# C2: $1,554.00 -> $1,558.00
$wsStands.Range("C2").Value = "'$1,558.00"
$wsStands.Range("C2").Style = "Normal"

# New row 6: Rogue SM-2.5 Monster Squat Stand 2.0
$wsStands.Range("A6").Value = "Rogue SM-2.5 Monster Squat Stand 2.0"
$wsStands.Range("B6").Value = "Rogue Fitness"
$wsStands.Range("C6").Value = "'Price not available"
$wsStands.Range("C6").Style = "Normal"
$wsStands.Range("D6").Value = "USA"

$imageUrl = "https://assets.roguefitness.com/f_auto,q_auto,c_limit,w_1042,b_rgb:f8f8f8/catalog/Rigs%20and%20Racks/Squat%20Stands/S2SQUAT2-0/S2SQUAT2-0-H_bmnpsa.png"
$productUrl = "https://www.roguefitness.com/rogue-sm-2-5-monster-squat-stand-2-0"

$wsStands.Hyperlinks.Add($wsStands.Range("E6"), $imageUrl) | Out-Null
$wsStands.Hyperlinks.Add($wsStands.Range("F6"), $productUrl) | Out-Null

# Match the existing hyperlink-cell look (blue/underline) used by E2:F5
$wsStands.Range("E2").Copy() | Out-Null
$wsStands.Range("E6").PasteSpecial(-4122) | Out-Null
$wsStands.Range("F2").Copy() | Out-Null
$wsStands.Range("F6").PasteSpecial(-4122) | Out-Null

Write-Host "Gym prices updated"
